# Fruta / hortaliza, semanal
# Insert a new week's worth of Kiwi "Gold" variety price records at rows 892-897
# (pushing the existing Hayward records down by 6 rows) in the
# "Femacal de La Calera - Kiwi" consolidated sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows before row 892, shifting everything below down.
$ws.Rows("892:897").Insert(-4121)

# Common (constant) values shared by every row in this block.
$mercadoId = 3
$mercado   = "Femacal de La Calera"
$region    = "Coquimbo"
$codreg    = 5
$tipo      = "Fruta"
$productoId = 100101
$producto   = "Berries"
$categoriaId = 100101007
$categoria   = "Kiwi"
$unidad    = "`$/bandeja 10 kilos"
$origen    = "Regi$([char]0x00F3)n de O'Higgins"
$kgUnidad  = 10

# New data for the "Gold" kiwi variety, week of 2023-03-23 (serial 45008).
$newRows = @(
    @{ Row = 892; Fecha = 45008; Calidad = "Especial"; Volumen = 56; Precio = 8000  },
    @{ Row = 893; Fecha = 45008; Calidad = "Primera";  Volumen = 67; Precio = 7000  },
    @{ Row = 894; Fecha = 45008; Calidad = "Segunda";  Volumen = 60; Precio = 6000  },
    @{ Row = 895; Fecha = 45008; Calidad = "Especial"; Volumen = 65; Precio = 10000 },
    @{ Row = 896; Fecha = 45008; Calidad = "Primera";  Volumen = 67; Precio = 8000  },
    @{ Row = 897; Fecha = 45008; Calidad = "Segunda";  Volumen = 60; Precio = 7000  }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.Fecha
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = "Gold"
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.Precio
    $ws.Cells.Item($row, 15).Value = $r.Precio
    $ws.Cells.Item($row, 16).Value = $r.Precio
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.Precio / 10
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
